# Actualización de horarios Línea 141 - 17/01/2026
# Nueva hora de scrapeo: 02:53:19

$wb = $excel.ActiveWorkbook

$oldTime = "02:41:48"
$newTime = "02:53:19"

# --- Hoja 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 4"

$ws1.Range("A6").Value = $newTime
$ws1.Range("D6").Value = 5

$ws1.Range("A7").Value = $newTime
$ws1.Range("D7").Value = 55

$ws1.Range("A8").Value = $newTime
$ws1.Range("D8").Value = 68

$ws1.Cells.Item(9, 1).Value = $newTime
$ws1.Cells.Item(9, 2).Value = "04:46"
$ws1.Cells.Item(9, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(9, 4).Value = 113
$ws1.Cells.Item(9, 5).Value = "LP1912"

# --- Hoja 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 2"

$ws2.Range("A6").Value = $newTime
$ws2.Range("D6").Value = 5

$ws2.Cells.Item(7, 1).Value = $newTime
$ws2.Cells.Item(7, 2).Value = "04:46"
$ws2.Cells.Item(7, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(7, 4).Value = 113
$ws2.Cells.Item(7, 5).Value = "LP1912"

# --- Hoja 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
